# Add newly-scraped Letterboxd entries (rows 55-65) to Sheet1, matching the
# "Add files via upload" commit: 11 new title/year/link rows appended below
# the existing 54, one of which (León, row 57) carries a real hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Title (A), Year (B), Link (D) - column C ("Blank") is never populated.
$titles = @(
    "Seret Efes",
    "Bottom Feeders",
    "Le√≥n",
    "Anyone Else But Me",
    "The Conspirators",
    "Night of the Coconut",
    "Mommy, I'm Scared",
    "Empire M",
    "The Book of Sun",
    "The Last Guest",
    "Cremaster 2"
)
$years = @(2015, 1999, 2013, 2024, 1972, 2022, 2004, 1972, 2020, 2019, 1999)
$links = @(
    "https://letterboxd.com/film/seret-efes/",
    "https://letterboxd.com/film/bottom-feeders-1999/",
    "https://letterboxd.com/film/leon/",
    "https://letterboxd.com/film/anyone-else-but-me/",
    "https://letterboxd.com/film/the-conspirators-1972/",
    "https://letterboxd.com/film/night-of-the-coconut-2022/",
    "https://letterboxd.com/film/mommy-im-scared/",
    "https://letterboxd.com/film/empire-m/",
    "https://letterboxd.com/film/the-book-of-sun/",
    "https://letterboxd.com/film/the-last-guest/",
    "https://letterboxd.com/film/cremaster-2/"
)

$startRow = 55
for ($i = 0; $i -lt $titles.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $titles[$i]
    $ws.Cells.Item($r, 2).Value = $years[$i]
    $ws.Cells.Item($r, 4).Value = $links[$i]
}

# Only the León row (57) got turned into a real clickable hyperlink in the
# source workbook; the rest stay as plain URL text, matching rows 55/56/58-65.
$linkCell = $ws.Cells.Item(57, 4)
$linkCell.Value = $links[2]
$ws.Hyperlinks.Add($linkCell, $links[2])
$linkCell.Style = "Hyperlink"

# Restore the selection the author left the workbook on (bottom of the list).
$ws.Range("L59").Select()
